# Generate Report for Handoff
# Replaces the placeholder png/md sample rows with the real calleeMd1/calleeMd2/
# callerMd1/callerMd2 localization handoff rows, and adds a 4th data row
# (row 5) to each of the three worksheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Base urls (pattern copied from the existing hyperlinks already present
# in the workbook; only the file name portion changes per row).
# ---------------------------------------------------------------------
$e2eBase = "https://github.com/OpenLocalizationTest/oltest/blob/da8ae7ad2a8327b18e1acbd967227cd1b325ad84/e2e"
$zhBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d3a4c9db020b8566fe698728ee2bc58c0aee08b2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eebc3ed1d50d6cc6401578fc85b1eeb3e6268c2f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$ws = $wb.Worksheets.Item("Overview")

# Clear every existing hyperlink on the sheet so we can rebuild them (and
# their underlying relationships) from scratch in a deterministic order.
$ws.Hyperlinks.Delete()

$overviewFiles = @("calleeMd1.md", "calleeMd2.md", "callerMd1.md", "callerMd2.md")

for ($i = 0; $i -lt $overviewFiles.Length; $i++) {
    $row  = $i + 2
    $file = $overviewFiles[$i]

    $ws.Hyperlinks.Add($ws.Range("A$row"), "$e2eBase/$file", [Type]::Missing, [Type]::Missing, $file)
    $ws.Range("B$row").Value2 = "Ready for handoff"
    $ws.Range("C$row").Value2 = "Ready for handoff"
    $ws.Range("D$row").Value2 = "2016-03-23 23:08:10"
    $ws.Range("D$row").NumberFormat = $dateFormat
}

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()

$zhRows = @(
    @{ File = "calleeMd1.md"; Xlf = "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf"; Dep = $null;                                    ReferencedBy = "e2e\callerMd2.md,`ne2e\callerMd1.md" },
    @{ File = "calleeMd2.md"; Xlf = "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf";  Dep = $null;                                    ReferencedBy = "e2e\callerMd1.md" },
    @{ File = "callerMd1.md"; Xlf = "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf";  Dep = "e2e\calleeMd1.md,`ne2e\calleeMd2.md";    ReferencedBy = $null },
    @{ File = "callerMd2.md"; Xlf = "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf";  Dep = "e2e\calleeMd1.md";                       ReferencedBy = $null }
)

for ($i = 0; $i -lt $zhRows.Length; $i++) {
    $row  = $i + 2
    $info = $zhRows[$i]

    $ws.Hyperlinks.Add($ws.Range("A$row"), "$e2eBase/$($info.File)", [Type]::Missing, [Type]::Missing, $info.File)
    $ws.Range("B$row").Value2 = ".md"
    $ws.Range("C$row").Value2 = "Ready for handoff"
    $ws.Hyperlinks.Add($ws.Range("D$row"), "$zhBase/$($info.Xlf)", [Type]::Missing, [Type]::Missing, $info.Xlf)
    $ws.Range("E$row").Value2 = "2016-03-23 23:08:01"
    $ws.Range("H$row").Value2 = "0001-01-01 00:00:00"

    if ($info.Dep -ne $null) {
        $ws.Range("I$row").Value2 = $info.Dep
    }

    $ws.Range("J$row").Value2 = "Include"

    if ($info.ReferencedBy -ne $null) {
        $ws.Range("K$row").Value2 = $info.ReferencedBy
    }
}

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()

$deRows = @(
    @{ File = "calleeMd1.md"; Xlf = "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf"; Dep = $null;                                    ReferencedBy = "e2e\callerMd2.md,`ne2e\callerMd1.md" },
    @{ File = "calleeMd2.md"; Xlf = "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf";  Dep = $null;                                    ReferencedBy = "e2e\callerMd1.md" },
    @{ File = "callerMd1.md"; Xlf = "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf";  Dep = "e2e\calleeMd1.md,`ne2e\calleeMd2.md";    ReferencedBy = $null },
    @{ File = "callerMd2.md"; Xlf = "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf";  Dep = "e2e\calleeMd1.md";                       ReferencedBy = $null }
)

for ($i = 0; $i -lt $deRows.Length; $i++) {
    $row  = $i + 2
    $info = $deRows[$i]

    $ws.Hyperlinks.Add($ws.Range("A$row"), "$e2eBase/$($info.File)", [Type]::Missing, [Type]::Missing, $info.File)
    $ws.Range("B$row").Value2 = ".md"
    $ws.Range("C$row").Value2 = "Ready for handoff"
    $ws.Hyperlinks.Add($ws.Range("D$row"), "$deBase/$($info.Xlf)", [Type]::Missing, [Type]::Missing, $info.Xlf)
    $ws.Range("E$row").Value2 = "2016-03-23 23:08:10"
    $ws.Range("H$row").Value2 = "0001-01-01 00:00:00"

    if ($info.Dep -ne $null) {
        $ws.Range("I$row").Value2 = $info.Dep
    }

    $ws.Range("J$row").Value2 = "Include"

    if ($info.ReferencedBy -ne $null) {
        $ws.Range("K$row").Value2 = $info.ReferencedBy
    }
}

Write-Host "Report generated for handoff."
